# feat: add 2022-Q4 data
#
# Inserts a new "2022-Q4" sheet (holding the latest quarterly fund-holding
# snapshot) right after the "总计" (summary) sheet and before "2022-Q3",
# then updates the "总计" sheet with a new leading row for 2022-Q4 while
# shifting the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet by duplicating "2022-Q3" (same
#    column layout/styling as every other quarterly sheet) and placing
#    the copy immediately before it, then rename + fill in the new data.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

$q4.Range("B2").Value = "'004352"
$q4.Range("C2").Value = "北信瑞丰研究精选股票"
$q4.Range("D2").Value = "'0.01"
$q4.Range("E2").Value = "'92.09"
$q4.Range("F2").Value = "'1.24"
$q4.Range("G2").Value = "'0.0001"
# H2 (仓位排名) keeps the value inherited from the copied sheet: 8

# ---------------------------------------------------------------------
# 2) Update "总计" summary sheet: shift rows 2-5 down to rows 3-6
#    (bottom-up, cell by cell, so per-cell styling is preserved) and
#    write the new 2022-Q4 summary row into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

for ($row = 5; $row -ge 2; $row--) {
    $destRow = $row + 1
    foreach ($col in @("A", "B", "C", "D")) {
        $src = $col + $row
        $dst = $col + $destRow
        $summary.Range($src).Copy($summary.Range($dst))
    }
}

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0

# Re-number the running index column (A) sequentially 0..4 top to bottom,
# matching the target layout.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
